$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (15) -- shifts Length(mm), Weight(kg),
# Precocity, Tissue Sample, Comments one column to the right.
$ws.Range("O1").EntireColumn.Insert()

# Set the new header cell value and comment.
$ws.Range("O3").Value = "Lifestage"
$ws.Range("O3").AddComment("Use full name`nEg. Fry/Parr/Smolt")

Write-Output "done"
